# "minor update before release"
# Add a new "ITEM CATEGORY" column (E) to the import template header row,
# matching the look (bold header style) and sizing of the existing columns,
# and leave the cursor on the new column as the last-used selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, styled like the other bold header cells (B1:D1)
$ws.Range("E1").Value = "ITEM CATEGORY"
$ws.Range("E1").Font.Bold = $true

# Size the new column similarly to the others (best-fit width for the header text)
$ws.Columns.Item(5).ColumnWidth = 17.5

# Leave the selection on the cell below the new header, as last saved
$ws.Range("E4").Select()

# Restore the workbook window size/position recorded in the file
$win = $wb.Windows.Item(1)
$win.Left = 450
$win.Top = 1620
$win.Width = 21600
$win.Height = 11385
